$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a new "2022-Q4" sheet, positioned right before "2022-Q3", by
#    duplicating the "2022-Q3" sheet (so it inherits the exact same layout,
#    column widths, and cell styles) and then overwriting its data values.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4.Range("D2").Value = "'0.23"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'94.47"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'2.52"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.0058"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 5

# Restore the originally-selected tab ("2020-Q4") since copying a sheet
# makes the new copy the active tab.
$wb.Worksheets.Item("2020-Q4").Activate()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    "2022-Q4" and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2020-Q4"
$total.Range("C8").Value = 2
$total.Range("D8").Value = 0.05

$total.Range("B7").Value = "2021-Q3"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 0.08

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 0.06

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.16

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01
